$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.134.69"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "3.093.64"
$ws.Range("E3").Value = "  +0.90%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.092.00"
$ws.Range("E8").Value = "  +0.99%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.79%  "
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "3.594.19"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "64.164.64"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "3.090.59"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "483.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.01"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.676"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  +4.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +10.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.02"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.20%  "
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.34"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "455.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +16.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0409"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0822"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").Value = "2.985.57"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.27"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  -3.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.99"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  +2.64%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.14"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.24%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.56"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("D50").Value = "0.0₃0517"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("E51").Value = "  +0.56%  "
